$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = "50, 61, 55"
$ws.Range("D4").Value = "2, AVI, 8"
$ws.Range("N4").Value = "10, 70, 46, 55"
$ws.Range("R4").Value = "68, 5, 14"
$ws.Range("D5").Value = "AII/A, AXIII, AVI, AVII"
$ws.Range("G5").Value = "AV, AIII/6, AIII/8"
$ws.Range("I5").Value = "AVIII, AIII, 24"
$ws.Range("J5").Value = "AIII/4, AIII/6, AIII/3, 7"
$ws.Range("M5").Value = "AI, AII"
$ws.Range("N5").Value = "28, 29, 61, 72"
$ws.Range("O5").Value = "57, 68, 23, 24"
$ws.Range("Q5").Value = "AV, 2, 21, 5"
$ws.Range("R5").Value = "25, AIII/App2, 12, 13"
$ws.Range("S5").Value = "AIII, 11, 4, 6"
$ws.Range("T5").Value = "74, 115, 43, 62"
$ws.Range("J6").Value = "AI/PO, 20, AIII/3, 7"
$ws.Range("N7").Value = "10, 71, 55, 72"
$ws.Range("P9").Value = "AIII, 27, 8"
$ws.Range("S9").Value = "11, 4, 6"
$ws.Range("O11").Value = "9, 28, 41, 33"
$ws.Range("N12").Value = "55, 50, 22, 2"
$ws.Range("S12").Value = "9, 2, 11"
$ws.Range("J14").Value = "AIII/4, 33, AIII/3, 7"
$ws.Range("N14").Value = "10, 28, 13, 55"
$ws.Range("P14").Value = "35, AII/II, AII/III, 24"
$ws.Range("Q14").Value = "18, 53, 42"
$ws.Range("K15").Value = "73, 90, 92, 108"
$ws.Range("K17").Value = "AXI/A, 52, AX"
$ws.Range("B18").Value = "26, 19"
$ws.Range("G18").Value = "AIII/11, 14, AIII/1, AIII/2"
$ws.Range("M18").Value = "AII, 11, 13"
$ws.Range("P18").Value = "AIII, 15, AII/V"
$ws.Range("T18").Value = "115, 68, 45, 15"
$ws.Range("R19").Value = "AIII/App2, AIV"
$ws.Range("T19").Value = "68, 76, 62"
$ws.Range("O22").Value = "33, 11, 49, 8"
$ws.Range("O23").Value = "11, 45, 38, 63"
$ws.Range("A24").Value = "AIV/A, AV/B, 6"
$ws.Range("B24").Value = "26, 36, AIV"
$ws.Range("C24").Value = "30, 7, 15"
$ws.Range("D24").Value = "AII/A, AIV/A, AXIII, 8"
$ws.Range("E24").Value = "21, AXI, 13, 6"
$ws.Range("F24").Value = "10, AVI, 5"
$ws.Range("G24").Value = "AIV, AV, 6, 8"
$ws.Range("H24").Value = "26, 11, 13, AIII"
$ws.Range("I24").Value = "8, 24, AV"
$ws.Range("J24").Value = "9, AIII/1, 7, 23"
$ws.Range("K24").Value = "92, AVI/A, 31"
$ws.Range("L24").Value = "28, 52, AVI/A, 87"
$ws.Range("M24").Value = "AII, 4, 7"
$ws.Range("N24").Value = "10, 29, AV, 72"
$ws.Range("O24").Value = "9, 25, 24"
$ws.Range("P24").Value = "17, AIII, 35, 19"
$ws.Range("Q24").Value = "AV, 47, 48"
$ws.Range("R24").Value = "25, 26, 12"
$ws.Range("S24").Value = "4, 11, AIII, 6"
$ws.Range("T24").Value = "74, 45, 134, 56"
$ws.Range("T25").Value = "43, 85, 62"
